$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2763.7273
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 2763.7273
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 8291.1819
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -8515.1819
$ws.Range("H70").Value = 50750
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 50750
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 152250
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -152790
$ws.Range("H73").Value = 50750
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 50750
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 152250
$ws.Range("M73").Value = -154122
$ws.Range("H74").Value = 5423.778
$ws.Range("I74").Value = 5773.4287
$ws.Range("K74").Value = 5773.4287
$ws.Range("M74").Value = -4837.4287
$ws.Range("H77").Value = 5423.778
$ws.Range("I77").Value = 5773.4287
$ws.Range("K77").Value = 28867.1435
$ws.Range("M77").Value = -24187.1435
$ws.Range("H80").Value = 1760.5
$ws.Range("I80").Value = 2800.6667
$ws.Range("J80").Value = 720.3333
$ws.Range("K80").Value = 8402.000100000001
$ws.Range("L80").Value = 2160.9999
$ws.Range("M80").Value = -7404.000100000001
$ws.Range("N80").Value = -4156.9999
$ws.Range("H83").Value = 1760.5
$ws.Range("I83").Value = 2800.6667
$ws.Range("J83").Value = 720.3333
$ws.Range("K83").Value = 25206.0003
$ws.Range("L83").Value = 6482.9997
$ws.Range("M83").Value = -20214.0003
$ws.Range("N83").Value = -16466.9997
$ws.Range("H121").Value = 999
$ws.Range("J121").Value = 999
$ws.Range("L121").Value = 2997
$ws.Range("N121").Value = -6491
$ws.Range("H129").Value = 912.8333
$ws.Range("I129").Value = 1149.75
$ws.Range("J129").Value = 887.8946999999999
$ws.Range("K129").Value = 3449.25
$ws.Range("L129").Value = 2663.6841
$ws.Range("M129").Value = 1550.75
$ws.Range("N129").Value = -12663.6841
$ws.Range("H132").Value = 1024.5143
$ws.Range("I132").Value = 1029.6333
$ws.Range("K132").Value = 3088.8999
$ws.Range("M132").Value = -558.8998999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 1498
$ws.Range("J9").Value = 1498
$ws.Range("L9").Value = 1498
$ws.Range("N9").Value = -1838
$ws.Range("H20").Value = 1498
$ws.Range("J20").Value = 1498
$ws.Range("L20").Value = 1498
$ws.Range("N20").Value = -2038
$ws.Range("H32").Value = 3844.7454
$ws.Range("I32").Value = 2065.4187
$ws.Range("K32").Value = 2065.4187
$ws.Range("M32").Value = -1778.4187
$ws.Range("H45").Value = 1540.25
$ws.Range("I45").Value = 1080.2858
$ws.Range("J45").Value = 1898
$ws.Range("K45").Value = 1080.2858
$ws.Range("L45").Value = 1898
$ws.Range("M45").Value = -703.2858000000001
$ws.Range("N45").Value = -2652
$ws.Range("H74").Value = 1194.5
$ws.Range("I74").Value = 1194.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1194.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -320.5
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1194.5
$ws.Range("I77").Value = 1194.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5972.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1604.5
$ws.Range("N77").ClearContents()
$ws.Range("H102").Value = 1983
$ws.Range("I102").Value = 1461.3
$ws.Range("K102").Value = 1461.3
$ws.Range("M102").Value = 160.7
$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 113.55556
$ws.Range("J7").Value = 65.75
$ws.Range("L7").Value = 65.75
$ws.Range("N7").Value = -291.75
$ws.Range("H22").Value = 866.3333
$ws.Range("I22").Value = 299.33334
$ws.Range("K22").Value = 299.33334
$ws.Range("M22").Value = 50.66665999999998
$ws.Range("H31").Value = 3143.5417
$ws.Range("I31").Value = 1859.8182
$ws.Range("J31").Value = 4229.769
$ws.Range("K31").Value = 1859.8182
$ws.Range("L31").Value = 4229.769
$ws.Range("M31").Value = -1564.8182
$ws.Range("N31").Value = -4819.769
$ws.Range("H34").Value = 3143.5417
$ws.Range("I34").Value = 1859.8182
$ws.Range("J34").Value = 4229.769
$ws.Range("K34").Value = 1859.8182
$ws.Range("L34").Value = 4229.769
$ws.Range("M34").Value = -1657.8182
$ws.Range("N34").Value = -4633.769
$ws.Range("H58").Value = 1403819
$ws.Range("I58").Value = 2558583.8
$ws.Range("J58").Value = 1604.5714
$ws.Range("K58").Value = 2558583.8
$ws.Range("L58").Value = 1604.5714
$ws.Range("M58").Value = -2558380.8
$ws.Range("N58").Value = -2010.5714
$ws.Range("H74").Value = 25399.834
$ws.Range("J74").Value = 25399.834
$ws.Range("L74").Value = 25399.834
$ws.Range("N74").Value = -27147.834
$ws.Range("H77").Value = 25399.834
$ws.Range("J77").Value = 25399.834
$ws.Range("L77").Value = 76199.50199999999
$ws.Range("N77").Value = -84935.50199999999
$ws.Range("H105").Value = 1272.2858
$ws.Range("I105").Value = 979
$ws.Range("K105").Value = 979
$ws.Range("M105").Value = 768
$ws.Range("H132").Value = 2004.0555
$ws.Range("I132").Value = 1135.4615
$ws.Range("J132").Value = 4262.4
$ws.Range("K132").Value = 3406.3845
$ws.Range("L132").Value = 12787.2
$ws.Range("M132").Value = -876.3844999999997
$ws.Range("N132").Value = -17847.2
$ws.Range("H134").Value = 1885.3529
$ws.Range("I134").Value = 1378.25
$ws.Range("K134").Value = 4134.75
$ws.Range("M134").Value = -1599.75
$ws.Range("H136").Value = 1403819
$ws.Range("I136").Value = 2558583.8
$ws.Range("J136").Value = 1604.5714
$ws.Range("K136").Value = 7675751.399999999
$ws.Range("L136").Value = 4813.7142
$ws.Range("M136").Value = -7673201.399999999
$ws.Range("N136").Value = -9913.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 500
$ws.Range("J2").Value = 100
$ws.Range("L2").Value = 600
$ws.Range("N2").Value = -826
$ws.Range("H22").Value = 3000
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 3000
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H34").Value = 200.4
$ws.Range("I34").Value = 200.4
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 601.2
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -517.2
$ws.Range("N34").ClearContents()
$ws.Range("H39").Value = 3000
$ws.Range("J39").Value = 3000
$ws.Range("L39").Value = 9000
$ws.Range("N39").Value = -9588
$ws.Range("H55").Value = 50252
$ws.Range("I55").Value = 50252
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 150756
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -150579
$ws.Range("N55").ClearContents()
$ws.Range("H68").Value = 677.5
$ws.Range("I68").Value = 555
$ws.Range("K68").Value = 1665
$ws.Range("M68").Value = -854
$ws.Range("H71").Value = 677.5
$ws.Range("I71").Value = 555
$ws.Range("K71").Value = 4995
$ws.Range("M71").Value = -939
$ws.Range("H131").Value = 782.7041
$ws.Range("J131").Value = 807.1867999999999
$ws.Range("L131").Value = 2421.5604
$ws.Range("N131").Value = -12501.5604

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J5").Value = 9166.666999999999
$ws.Range("L5").Value = 9166.666999999999
$ws.Range("N5").Value = -9390.666999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 349916.66
$ws.Range("J2").Value = 114071.43
$ws.Range("L2").Value = 114071.43
$ws.Range("N2").Value = -114295.43
$ws.Range("H55").Value = 373.06668
$ws.Range("I55").Value = 319.42856
$ws.Range("K55").Value = 319.42856
$ws.Range("M55").Value = -146.42856
$ws.Range("H64").Value = 512399.5
$ws.Range("J64").Value = 24800
$ws.Range("L64").Value = 24800
$ws.Range("N64").Value = -25250
$ws.Range("H67").Value = 512399.5
$ws.Range("J67").Value = 24800
$ws.Range("L67").Value = 24800
$ws.Range("N67").Value = -26360
$ws.Range("H136").Value = 4730.231
$ws.Range("I136").Value = 2199.75
$ws.Range("J136").Value = 5854.8887
$ws.Range("K136").Value = 6599.25
$ws.Range("L136").Value = 17564.6661
$ws.Range("M136").Value = -4049.25
$ws.Range("N136").Value = -22664.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 12124.75
$ws.Range("J41").Value = 12124.75
$ws.Range("L41").Value = 12124.75
$ws.Range("N41").Value = -12904.75
$ws.Range("H45").Value = 13428.167
$ws.Range("J45").Value = 12000
$ws.Range("L45").Value = 12000
$ws.Range("N45").Value = -12982
